$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.870.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.33%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.274.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.00%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.76%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.34%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('E8').Value = '  -3.28%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.266.07'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.93%  '

$ws.Range('E10').Value = '  -9.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.585'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.39%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.30'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000265'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.60'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.97%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '633.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.805.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.68%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.861.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.69%  '

$ws.Range('E19').Value = '  -3.43%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.277.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.89%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.33%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.903'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.73%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '107.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.88'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.32%  '

$ws.Range('E26').Value = '  -7.75%  '

$ws.Range('E27').Value = '  -7.61%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.26%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.18%  '

$ws.Range('E31').Value = '  -6.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.38%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.27%  '

$ws.Range('E34').Value = '  -4.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.12%  '

$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.694.36'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.75%  '

$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '521.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.35%  '

$ws.Range('E39').Value = '  -4.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0727'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.85%  '

$ws.Range('E41').Value = '  -2.75%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '32.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.76%  '

$ws.Range('B44').Value = 'CoreDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.28%  '

$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.337'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.98%  '

$ws.Range('E46').Value = '  -1.71%  '

$ws.Range('E47').Value = '  -6.61%  '

$ws.Range('E48').Value = '  -4.25%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.75%  '

$ws.Range('E50').Value = '  +0.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.67%  '
